$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("P5").Value = 0.3870967741935484
$ws.Range("P6").Value = 0.5355846774193548
$ws.Range("P7").Value = 0.07949583949697754
$ws.Range("P8").Value = 0.7384604291930315
